$d = $word.ActiveDocument

# Locate the last paragraph in the document ("To make sure we can't do
# up/down or anything related to text if we still don't have text on
# the meme yet") and append two new list items after it, matching the
# ListParagraph style / numbering (ilvl 0, numId 1) already used by the
# other bullet points.

$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()

$p1 = $d.Paragraphs.Last
$p1.Range.InsertAfter("Font-size/family/letter-spacing etc in the editor (everything)")

$p1.Range.InsertParagraphAfter()

$p2 = $d.Paragraphs.Last
$p2.Range.InsertAfter("Add flexbox helpers and use them in HTML")
